$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2,3,4 are being cyclically rotated:
#   new row2 = old row3
#   new row3 = old row4
#   new row4 = old row2
# Columns A,B,C,E,F,G,H,I,N,Q,R are identical across the three rows, so only
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) actually change.
# Note: use Value2 for reads (plain .Value getter is unreliable in this runtime).

$oldD2 = $ws.Range("D2").Value2
$oldJ2 = $ws.Range("J2").Value2
$oldK2 = $ws.Range("K2").Value2
$oldL2 = $ws.Range("L2").Value2
$oldM2 = $ws.Range("M2").Value2
$oldO2 = $ws.Range("O2").Value2
$oldP2 = $ws.Range("P2").Value2

$oldD3 = $ws.Range("D3").Value2
$oldJ3 = $ws.Range("J3").Value2
$oldK3 = $ws.Range("K3").Value2
$oldL3 = $ws.Range("L3").Value2
$oldM3 = $ws.Range("M3").Value2
$oldO3 = $ws.Range("O3").Value2
$oldP3 = $ws.Range("P3").Value2

$oldD4 = $ws.Range("D4").Value2
$oldJ4 = $ws.Range("J4").Value2
$oldK4 = $ws.Range("K4").Value2
$oldL4 = $ws.Range("L4").Value2
$oldM4 = $ws.Range("M4").Value2
$oldO4 = $ws.Range("O4").Value2
$oldP4 = $ws.Range("P4").Value2

# Row 2 <- old Row 3
$ws.Range("D2").Value2 = $oldD3
$ws.Range("J2").Value2 = $oldJ3
$ws.Range("K2").Value2 = $oldK3
$ws.Range("L2").Value2 = $oldL3
$ws.Range("M2").Value2 = $oldM3
$ws.Range("O2").Value2 = $oldO3
$ws.Range("P2").Value2 = $oldP3

# Row 3 <- old Row 4
$ws.Range("D3").Value2 = $oldD4
$ws.Range("J3").Value2 = $oldJ4
$ws.Range("K3").Value2 = $oldK4
$ws.Range("L3").Value2 = $oldL4
$ws.Range("M3").Value2 = $oldM4
$ws.Range("O3").Value2 = $oldO4
$ws.Range("P3").Value2 = $oldP4

# Row 4 <- old Row 2
$ws.Range("D4").Value2 = $oldD2
$ws.Range("J4").Value2 = $oldJ2
$ws.Range("K4").Value2 = $oldK2
$ws.Range("L4").Value2 = $oldL2
$ws.Range("M4").Value2 = $oldM2
$ws.Range("O4").Value2 = $oldO2
$ws.Range("P4").Value2 = $oldP2
